$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("out_vars")

# --- New data rows: raw and clean SSA data through June 11th ---

# Row 11 (2020-06-10)
$ws.Cells.Item(11, 1).Value = 9
$ws.Cells.Item(11, 2).Value = 43992
$ws.Cells.Item(11, 3).Value = 129184
$ws.Cells.Item(11, 4).Value = 186570
$ws.Cells.Item(11, 5).Value = 53608
$ws.Cells.Item(11, 6).Value = 15357
$ws.Cells.Item(11, 7).Value = 33.11

# Row 12 (2020-06-11)
$ws.Cells.Item(12, 1).Value = 10
$ws.Cells.Item(12, 2).Value = 43993
$ws.Cells.Item(12, 3).Value = 133974
$ws.Cells.Item(12, 4).Value = 191465
$ws.Cells.Item(12, 5).Value = 55700
$ws.Cells.Item(12, 6).Value = 15944
$ws.Cells.Item(12, 7).Value = 33.01

# Copy the style (border/bold) used by column A further down (A2:A10) onto the new index cells
$ws.Range("A10").Copy() | Out-Null
$ws.Range("A11:A12").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

# Update the date column's number format (applies to the whole Fecha column,
# old and new rows alike) from the short date to a full timestamp format
$ws.Range("B2:B12").NumberFormat = "yyyy\-mm\-dd\ hh:mm:ss"

# Selection ends up on B13, just below the last data row, as in the source workbook
$ws.Range("B13").Select() | Out-Null
